# Vietnamese translation pass for
# "Email 3 [TEMPLATE] Partner email - list of travel documents.docx"
#
# Strategy: most English sentences/phrases occur exactly once in the
# document's flattened text, so a simple Find/ReplaceOne on $d.Content is
# safe and order-independent. A few short phrases ("and ", " on ", " or ")
# are ambiguous (appear more than once with possibly different
# replacements), so those are handled with range-scoped, sequential finds
# that walk forward through the document. Comment bubbles are translated
# by setting Comments.Item(n).Range.Text directly.

$d = $word.ActiveDocument

function Replace-Unique($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 1) | Out-Null
}

# --- Simple, unique, order-independent replacements -----------------------

Replace-Unique "An email sent to partners in the target country who have RSVPed yes. We want them to submit their documents. It will be sent via customer.io" "Một email gửi đến các đối tác tại quốc gia mục tiêu đã xác nhận tham gia. Chúng tôi muốn họ gửi tài liệu của mình. It will be sent via customer.io"

Replace-Unique "Invited partners who RSVP yes" "Các đối tác được mời đã xác nhận tham gia"

Replace-Unique "Subject: " "Chủ đề: "

Replace-Unique " — take the next step" " — Bước tiếp theo"

Replace-Unique "Thank you for registering for " "Cảm ơn bạn đã đăng ký tham dự sự kiện "

Replace-Unique "Hi " "Xin chào "

Replace-Unique "We are excited for you to join us at " "Chúng tôi rất vui khi bạn có thể tham dự với chúng tôi tại sự kiện "

Replace-Unique "To confirm your registration, we would require you and one guest of your choice to provide us with:" "Để tiến hành đăng ký cho bạn, chúng tôi cần bạn và một khách mời của bạn cung cấp cho chúng tôi:"

Replace-Unique "A signed copy of the " "Một bản sao có chữ ký bộ "

Replace-Unique "Code of Conduct " "Quy tắc Ứng xử"

# "and " is ambiguous (occurs several times in the flattened text); handled
# below with a range scoped right after "Quy tắc Ứng xử".

Replace-Unique "Terms and Conditions" "Điều khoản và Điều kiện"

Replace-Unique " (1 set from each person)" " (mỗi người 1 bộ)"

Replace-Unique "A scanned copy of your international passports" "Bản scan hộ chiếu quốc tế của bạn"

Replace-Unique "Covid-19 vaccination certificates" "Giấy chứng nhận tiêm phòng Covid-19"

Replace-Unique "Send my details" "Gửi thông tin của tôi"

Replace-Unique "Your country manager will be in touch to confirm your booking or request any other relevant details. " "Giám đốc phụ trách tại quốc gia của bạn sẽ liên lạc để xác nhận hoặc hỏi thêm các thông tin liên quan khác nếu cần. "

Replace-Unique "Our event package offers you and your guest: " "Gói sự kiện chúng tôi cung cấp đến bạn và khách mời của bạn bao gồm: "

Replace-Unique "Flight tickets " "Vé máy bay "

Replace-Unique "Travel insurance " "Bảo hiểm du lịch "

Replace-Unique "Airport – Hotel – Airport transfer " "Đưa đón sân bay – khách sạn "

Replace-Unique "One hotel room for you and your guest / Two hotel rooms for you and your guest" "Một hoặc hai phòng khách sạn cho bạn và khách mời của bạn"

Replace-Unique "Check-in" "Nhận phòng"

Replace-Unique "Check-out" "Trả phòng"

# " on " occurs twice (Check-in / Check-out) but both map to the same
# Vietnamese text, so ReplaceAll is safe.
$d.Content.Find.Execute(" on ", $true, $false, $false, $false, $false, $true, 1, $false, " vào ngày ", 2) | Out-Null

Replace-Unique "Meals (Breakfast, lunch, and dinner)" "Các bữa ăn (Bữa sáng, bữa trưa và bữa tối)"

Replace-Unique "Sightseeing tour of " "Tour tham quan "

Replace-Unique "We will send you a confirmation letter before your departure date with the event agenda and information about your flights, transportation, and accommodation. " "Chúng tôi sẽ gửi thư xác nhận đến bạn trước ngày khởi hành với các thông tin chi tiết về chương trình sự kiện, chuyến bay, phương tiện di chuyển và chỗ ở của bạn. "

Replace-Unique "If you have any questions, please contact us via " "Nếu bạn cần hỗ trợ, vui lòng liên hệ với chúng tôi qua "

Replace-Unique "If you have any questions, please contact your country manager, " "Nếu bạn có bất kỳ thắc mắc nào, vui lòng liên hệ với giám đốc phụ trách quốc gia của bạn "

Replace-Unique ", at " ", qua email "

Replace-Unique "We look forward to seeing you soon." "Chúng tôi rất mong được gặp bạn."

# --- Ambiguous short phrases, resolved with a forward-walking range -------

# "and " — only the one right after "Quy tắc Ứng xử" (formerly "Code of
# Conduct ") needs replacing, becoming "và ".
$rng = $d.Content
$rng.Find.Execute("Quy tắc Ứng xử", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.Find.Execute("and ", $true, $false, $false, $false, $false, $true, 1, $false, "và ", 1) | Out-Null

# " or " — two occurrences with different Vietnamese replacements. Anchor
# on text that is never translated (the "live chat" hyperlink, and the
# "[EMAIL ADDRESS]" placeholder) so the lookup stays valid no matter what
# order the surrounding sentences were translated in.
#
# The run between the "live chat" and "WhatsApp" hyperlinks carries no
# rPr of its own (plain), so replace only the inner "or" (leaving the
# surrounding spaces untouched) to avoid the replacement inheriting the
# adjacent hyperlink's run formatting.
$h2 = $d.Hyperlinks.Item(2)  # "live chat"
$rngOr1 = $d.Range($h2.Range.End + 1, $h2.Range.End + 3)
$rngOr1.Text = "hoặc"

$rng3 = $d.Content
$rng3.Find.Execute("[EMAIL ADDRESS]", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng3.Collapse(0)
$rng3.Find.Execute(" or ", $true, $false, $false, $false, $false, $true, 1, $false, " hoặc số ", 1) | Out-Null

# --- Comment bubbles --------------------------------------------------

# Comments.Item(n) is ordered by position of commentRangeStart in the
# document body: 1="please check if..." 2="link to COC" 3="link to T&C"
# 4="please confirm these" 5="choose either one" (left untranslated).
$d.Comments.Item(1).Range.Text = "vui lòng kiểm tra xem đây có phải là tất cả các tài liệu cần thiết không"
$d.Comments.Item(2).Range.Text = "liên kết đến COC"
$d.Comments.Item(3).Range.Text = "liên kết đến T&C"
$d.Comments.Item(4).Range.Text = "vui lòng xác nhận những điều này"

Write-Host "Done."
